$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad) holds a date serial number for rows 2-72.
# It was updated from 45178 (2023-09-09) to 45179 (2023-09-10) on every row.
for ($row = 2; $row -le 72; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45178) {
        $cell.Value2 = 45179
    }
}
